$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Hsp90aa1"
$ws.Range("C2").Value = "Fgfr3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 107.928617
$ws.Range("H2").Value = 323.785851
$ws.Range("I2").Value = 0.2068777607879145
$ws.Range("J2").Value = 0.2068777607879145
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 2.619953333333334
$ws.Range("N2").Value = 7.85986
$ws.Range("O2").Value = 0.6100029538328192
$ws.Range("P2").Value = 0.6100029538328192
$ws.Range("Q2").Value = 282.7679398712066
$ws.Range("R2").Value = 2544.91145884086
$ws.Range("S2").Value = 0.1261960451629472
$ws.Range("T2").Value = 0.1261960451629472

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Hsp90aa1"
$ws.Range("C3").Value = "Fgfr3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 107.928617
$ws.Range("H3").Value = 323.785851
$ws.Range("I3").Value = 0.2068777607879145
$ws.Range("J3").Value = 0.2068777607879145
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.5698483333333333
$ws.Range("N3").Value = 1.709545
$ws.Range("O3").Value = 0.1326776176306101
$ws.Range("P3").Value = 0.1326776176306101
$ws.Range("Q3").Value = 61.50294251642165
$ws.Range("R3").Value = 553.526482647795
$ws.Range("S3").Value = 0.02744804844209573
$ws.Range("T3").Value = 0.02744804844209573

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Hsp90aa1"
$ws.Range("C4").Value = "Fgfr3"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 107.928617
$ws.Range("H4").Value = 323.785851
$ws.Range("I4").Value = 0.2068777607879145
$ws.Range("J4").Value = 0.2068777607879145
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.105183
$ws.Range("N4").Value = 3.315549
$ws.Range("O4").Value = 0.2573194285365706
$ws.Range("P4").Value = 0.2573194285365706
$ws.Range("Q4").Value = 119.280872721911
$ws.Range("R4").Value = 1073.527854497199
$ws.Range("S4").Value = 0.05323366718287151
$ws.Range("T4").Value = 0.05323366718287151

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Hsp90aa1"
$ws.Range("C5").Value = "Fgfr3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 143.300008
$ws.Range("H5").Value = 429.900024
$ws.Range("I5").Value = 0.2746777045788536
$ws.Range("J5").Value = 0.2746777045788536
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 2.619953333333334
$ws.Range("N5").Value = 7.85986
$ws.Range("O5").Value = 0.6100029538328192
$ws.Range("P5").Value = 0.6100029538328192
$ws.Range("Q5").Value = 375.4393336262934
$ws.Range("R5").Value = 3378.95400263664
$ws.Range("S5").Value = 0.1675542111451192
$ws.Range("T5").Value = 0.1675542111451192

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Hsp90aa1"
$ws.Range("C6").Value = "Fgfr3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 143.300008
$ws.Range("H6").Value = 429.900024
$ws.Range("I6").Value = 0.2746777045788536
$ws.Range("J6").Value = 0.2746777045788536
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.5698483333333333
$ws.Range("N6").Value = 1.709545
$ws.Range("O6").Value = 0.1326776176306101
$ws.Range("P6").Value = 0.1326776176306101
$ws.Range("Q6").Value = 81.65927072545334
$ws.Range("R6").Value = 734.93343652908
$ws.Range("S6").Value = 0.03644358345976681
$ws.Range("T6").Value = 0.03644358345976681

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Hsp90aa1"
$ws.Range("C7").Value = "Fgfr3"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 143.300008
$ws.Range("H7").Value = 429.900024
$ws.Range("I7").Value = 0.2746777045788536
$ws.Range("J7").Value = 0.2746777045788536
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.105183
$ws.Range("N7").Value = 3.315549
$ws.Range("O7").Value = 0.2573194285365706
$ws.Range("P7").Value = 0.2573194285365706
$ws.Range("Q7").Value = 158.372732741464
$ws.Range("R7").Value = 1425.354594673176
$ws.Range("S7").Value = 0.07067990997396759
$ws.Range("T7").Value = 0.07067990997396759

$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Hsp90aa1"
$ws.Range("C8").Value = "Fgfr3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 134.5459086666667
$ws.Range("H8").Value = 403.637726
$ws.Range("I8").Value = 0.2578978317505474
$ws.Range("J8").Value = 0.2578978317505473
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 2.619953333333334
$ws.Range("N8").Value = 7.85986
$ws.Range("O8").Value = 0.6100029538328192
$ws.Range("P8").Value = 0.6100029538328192
$ws.Range("Q8").Value = 352.5040018975956
$ws.Range("R8").Value = 3172.536017078361
$ws.Range("S8").Value = 0.1573184391549133
$ws.Range("T8").Value = 0.1573184391549133

$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Hsp90aa1"
$ws.Range("C9").Value = "Fgfr3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 134.5459086666667
$ws.Range("H9").Value = 403.637726
$ws.Range("I9").Value = 0.2578978317505474
$ws.Range("J9").Value = 0.2578978317505473
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.5698483333333333
$ws.Range("N9").Value = 1.709545
$ws.Range("O9").Value = 0.1326776176306101
$ws.Range("P9").Value = 0.1326776176306101
$ws.Range("Q9").Value = 76.6707618105189
$ws.Range("R9").Value = 690.03685629467
$ws.Range("S9").Value = 0.03421726990876253
$ws.Range("T9").Value = 0.03421726990876253

$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Hsp90aa1"
$ws.Range("C10").Value = "Fgfr3"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 134.5459086666667
$ws.Range("H10").Value = 403.637726
$ws.Range("I10").Value = 0.2578978317505474
$ws.Range("J10").Value = 0.2578978317505473
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.105183
$ws.Range("N10").Value = 3.315549
$ws.Range("O10").Value = 0.2573194285365706
$ws.Range("P10").Value = 0.2573194285365706
$ws.Range("Q10").Value = 148.6978509779527
$ws.Range("R10").Value = 1338.280658801574
$ws.Range("S10").Value = 0.06636212268687149
$ws.Range("T10").Value = 0.06636212268687147

$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Hsp90aa1"
$ws.Range("C11").Value = "Fgfr3"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 135.927831
$ws.Range("H11").Value = 407.783493
$ws.Range("I11").Value = 0.2605467028826847
$ws.Range("J11").Value = 0.2605467028826847
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 2.619953333333334
$ws.Range("N11").Value = 7.85986
$ws.Range("O11").Value = 0.6100029538328192
$ws.Range("P11").Value = 0.6100029538328192
$ws.Range("Q11").Value = 356.12457392122
$ws.Range("R11").Value = 3205.12116529098
$ws.Range("S11").Value = 0.1589342583698396
$ws.Range("T11").Value = 0.1589342583698396

$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Hsp90aa1"
$ws.Range("C12").Value = "Fgfr3"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 135.927831
$ws.Range("H12").Value = 407.783493
$ws.Range("I12").Value = 0.2605467028826847
$ws.Range("J12").Value = 0.2605467028826847
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.5698483333333333
$ws.Range("N12").Value = 1.709545
$ws.Range("O12").Value = 0.1326776176306101
$ws.Range("P12").Value = 0.1326776176306101
$ws.Range("Q12").Value = 77.458247948965
$ws.Range("R12").Value = 697.1242315406849
$ws.Range("S12").Value = 0.034568715819985
$ws.Range("T12").Value = 0.034568715819985

$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Hsp90aa1"
$ws.Range("C13").Value = "Fgfr3"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 135.927831
$ws.Range("H13").Value = 407.783493
$ws.Range("I13").Value = 0.2605467028826847
$ws.Range("J13").Value = 0.2605467028826847
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.105183
$ws.Range("N13").Value = 3.315549
$ws.Range("O13").Value = 0.2573194285365706
$ws.Range("P13").Value = 0.2573194285365706
$ws.Range("Q13").Value = 150.225128048073
$ws.Range("R13").Value = 1352.026152432657
$ws.Range("S13").Value = 0.06704372869286009
$ws.Range("T13").Value = 0.06704372869286009
